# Updated for new functionality. Improved NMR data processing, overhauled
# MATLAB and dFBA code. Refresh the fitted coefficient/bound tables on the
# cf / lb / ub sheets and restore the UI selection state.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("cf")
$ws2 = $wb.Worksheets.Item("lb")
$ws3 = $wb.Worksheets.Item("ub")

# --- sheet1 (cf) value updates ---
$ws1.Range("B2").Value = 26.975441715296807
$ws1.Range("C2").Value = -0.13047633024004807
$ws1.Range("D2").Value = 16.431129540617491
$ws1.Range("E2").Value = 0.52455828470319477
$ws1.Range("B3").Value = 6.7687705930713369
$ws1.Range("C3").Value = 0.20866762468764199
$ws1.Range("D3").Value = 23.572524809817249
$ws1.Range("E3").Value = 0
$ws1.Range("B4").Value = 12.776293567874589
$ws1.Range("C4").Value = 0.2129993785622378
$ws1.Range("D4").Value = 15.123946984450281
$ws1.Range("E4").Value = 0
$ws1.Range("B5").Value = 14.00417375517492
$ws1.Range("C5").Value = 0.2440505078056821
$ws1.Range("D5").Value = 13.238862551799018
$ws1.Range("E5").Value = 0
$ws1.Range("B6").Value = 3.9928419306383338
$ws1.Range("C6").Value = 0.10736135248046683
$ws1.Range("D6").Value = 26.813621407055177
$ws1.Range("E6").Value = 0
$ws1.Range("B7").Value = 6.1815410660444092
$ws1.Range("C7").Value = 0.17324500460397269
$ws1.Range("D7").Value = 24.402508713684416
$ws1.Range("E7").Value = 0
$ws1.Range("B8").Value = 3.2128258949005866
$ws1.Range("C8").Value = 0.11125748733312468
$ws1.Range("D8").Value = 27.008809298244145
$ws1.Range("E8").Value = 0
$ws1.Range("B9").Value = 4.0045904723688448
$ws1.Range("C9").Value = -0.45735078498372128
$ws1.Range("D9").Value = 11.216599095619713
$ws1.Range("E9").Value = 2.9554095276311547
$ws1.Range("B10").Value = 3.9773745520001111
$ws1.Range("C10").Value = 0.64123545400628634
$ws1.Range("D10").Value = 11.875240026351921
$ws1.Range("E10").Value = 0
$ws1.Range("B11").Value = 7.63
$ws1.Range("C11").Value = -0.53374918531371707
$ws1.Range("D11").Value = 8.2882989136034873
$ws1.Range("E11").Value = 0
$ws1.Range("B12").Value = 2.1853930592703086
$ws1.Range("C12").Value = 0.71615649391571201
$ws1.Range("D12").Value = 4.6677561747542491
$ws1.Range("E12").Value = 0
$ws1.Range("B13").Value = 5.1832807094315205
$ws1.Range("C13").Value = 0.89761277885073454
$ws1.Range("D13").Value = 9.8755992348235733
$ws1.Range("E13").Value = 0

# --- sheet2 (lb) value updates ---
$ws2.Range("B2").Value = 26.530363328679737
$ws2.Range("C2").Value = -0.10748393003633491
$ws2.Range("D2").Value = 18.430296587256777
$ws2.Range("E2").Value = 0.96963667132026454
$ws2.Range("B3").Value = 6.3698442557035593
$ws2.Range("C3").Value = 0.2637715953076068
$ws2.Range("D3").Value = 25.047260317022705
$ws2.Range("E3").Value = 0
$ws2.Range("B4").Value = 12.186993683681795
$ws2.Range("C4").Value = 0.16830647718351471
$ws2.Range("D4").Value = 14.009889610938391
$ws2.Range("E4").Value = 0
$ws2.Range("B5").Value = 13.31039047557311
$ws2.Range("C5").Value = 0.18096725650509959
$ws2.Range("D5").Value = 12.044349782883744
$ws2.Range("E5").Value = 0
$ws2.Range("B6").Value = 3.2922075516615887
$ws2.Range("C6").Value = 0.055643578598840594
$ws2.Range("D6").Value = 21.232061932690634
$ws2.Range("E6").Value = 0
$ws2.Range("B7").Value = 5.7343715486349742
$ws2.Range("C7").Value = 0.12554883233522485
$ws2.Range("D7").Value = 22.5369398556875
$ws2.Range("E7").Value = 0
$ws2.Range("B8").Value = 2.5675857103219473
$ws2.Range("C8").Value = 0.03950081864463699
$ws2.Range("D8").Value = 19.812627165840262
$ws2.Range("E8").Value = 0
$ws2.Range("B9").Value = 3.9468306252662599
$ws2.Range("C9").Value = -0.52514868823439298
$ws2.Range("D9").Value = 10.902383734719804
$ws2.Range("E9").Value = 3.0131693747337405
$ws2.Range("B10").Value = 4.0162534813083939
$ws2.Range("C10").Value = 0.49791111319785514
$ws2.Range("D10").Value = 11.484727864283165
$ws2.Range("E10").Value = 0
$ws2.Range("B11").Value = 7.63
$ws2.Range("C11").Value = -0.64410115587848005
$ws2.Range("D11").Value = 7.8422898637383822
$ws2.Range("E11").Value = 0
$ws2.Range("B12").Value = 2.1369391201172645
$ws2.Range("C12").Value = 0.34521769080321613
$ws2.Range("D12").Value = 3.8404564183549121
$ws2.Range("E12").Value = 0
$ws2.Range("B13").Value = 5.2252688170940749
$ws2.Range("C13").Value = 0.70172750036117071
$ws2.Range("D13").Value = 9.596959097846355
$ws2.Range("E13").Value = 0

# --- sheet3 (ub) value updates ---
$ws3.Range("B2").Value = 27.558432267967394
$ws3.Range("C2").Value = -0.15346873044376122
$ws3.Range("D2").Value = 14.431962493978203
$ws3.Range("E2").Value = -0.058432267967392405
$ws3.Range("B3").Value = 7.271815100314253
$ws3.Range("C3").Value = 0.15356365406767719
$ws3.Range("D3").Value = 22.097789302611794
$ws3.Range("E3").Value = 0
$ws3.Range("B4").Value = 13.368155238456019
$ws3.Range("C4").Value = 0.25769227994096089
$ws3.Range("D4").Value = 16.238004357962172
$ws3.Range("E4").Value = 0
$ws3.Range("B5").Value = 14.701993655958612
$ws3.Range("C5").Value = 0.30713375910626461
$ws3.Range("D5").Value = 14.433375320714292
$ws3.Range("E5").Value = 0
$ws3.Range("B6").Value = 4.70753862803135
$ws3.Range("C6").Value = 0.15907912636209307
$ws3.Range("D6").Value = 32.395180881419719
$ws3.Range("E6").Value = 0
$ws3.Range("B7").Value = 6.6341111105320145
$ws3.Range("C7").Value = 0.22094117687272052
$ws3.Range("D7").Value = 26.268077571681332
$ws3.Range("E7").Value = 0
$ws3.Range("B8").Value = 3.8476799062968516
$ws3.Range("C8").Value = 0.18301415602161236
$ws3.Range("D8").Value = 34.204991430648029
$ws3.Range("E8").Value = 0
$ws3.Range("B9").Value = 4.0583597917256515
$ws3.Range("C9").Value = -0.38955288173304953
$ws3.Range("D9").Value = 11.530814456519623
$ws3.Range("E9").Value = 2.901640208274348
$ws3.Range("B10").Value = 3.9411817005454517
$ws3.Range("C10").Value = 0.7845597948147176
$ws3.Range("D10").Value = 12.265752188420677
$ws3.Range("E10").Value = 0
$ws3.Range("B11").Value = 7.63
$ws3.Range("C11").Value = -0.42339721474895409
$ws3.Range("D11").Value = 8.7343079634685914
$ws3.Range("E11").Value = 0
$ws3.Range("B12").Value = 2.2312095499167066
$ws3.Range("C12").Value = 1.087095297028208
$ws3.Range("D12").Value = 5.495055931153586
$ws3.Range("E12").Value = 0
$ws3.Range("B13").Value = 5.1435781016434285
$ws3.Range("C13").Value = 1.0934980573402984
$ws3.Range("D13").Value = 10.154239371800792
$ws3.Range("E13").Value = 0

# --- restore per-sheet selection state ---
$ws1.Activate()
$ws1.Range("B2").Select()

$ws2.Activate()
$ws2.Range("B2").Select()

# "ub" ends up the active/selected tab (activeTab=2 in workbook.xml)
$ws3.Activate()
$ws3.Range("B2").Select()
